$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 953.8333
$ws.Range("I2").Value = 599.5
$ws.Range("J2").Value = 1308.1666
$ws.Range("K2").Value = 599.5
$ws.Range("L2").Value = 1308.1666
$ws.Range("M2").Value = -486.5
$ws.Range("N2").Value = -1534.1666

# Row 9
$ws.Range("H9").Value = 320.33334
$ws.Range("I9").Value = 316.625
$ws.Range("K9").Value = 316.625
$ws.Range("M9").Value = -147.625

# Row 132
$ws.Range("H132").Value = 34485056
$ws.Range("I132").Value = 37039440
$ws.Range("K132").Value = 111118320
$ws.Range("M132").Value = -111115790

# Row 138
$ws.Range("H138").Value = 3009.49
$ws.Range("I138").Value = 1673.2122
$ws.Range("K138").Value = 5019.6366
$ws.Range("M138").Value = 120.3634000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8430.513999999999
$ws.Range("I32").Value = 6358.25
$ws.Range("K32").Value = 6358.25
$ws.Range("M32").Value = -6071.25

# Row 44
$ws.Range("H44").Value = 6947.75

# Row 45
$ws.Range("H45").Value = 6542413
$ws.Range("I45").Value = 8463536
$ws.Range("K45").Value = 8463536
$ws.Range("M45").Value = -8463159

# Row 88
$ws.Range("H88").Value = 1003.1429
$ws.Range("I88").Value = 958.3333
$ws.Range("J88").Value = 1036.75
$ws.Range("K88").Value = 958.3333
$ws.Range("L88").Value = 1036.75
$ws.Range("M88").Value = -552.3333
$ws.Range("N88").Value = -1848.75

# Row 91
$ws.Range("H91").Value = 1003.1429
$ws.Range("I91").Value = 958.3333
$ws.Range("J91").Value = 1036.75
$ws.Range("K91").Value = 958.3333
$ws.Range("L91").Value = 1036.75
$ws.Range("M91").Value = 445.6667
$ws.Range("N91").Value = -3844.75

# Row 110
$ws.Range("H110").Value = 896988.75
$ws.Range("I110").Value = 993022.4
$ws.Range("J110").Value = 675.3333
$ws.Range("K110").Value = 993022.4
$ws.Range("L110").Value = 675.3333
$ws.Range("M110").Value = -990977.4
$ws.Range("N110").Value = -4765.3333

# Row 132
$ws.Range("H132").Value = 5079.086
$ws.Range("I132").Value = 5288.2964
$ws.Range("J132").Value = 4373
$ws.Range("K132").Value = 15864.8892
$ws.Range("L132").Value = 13119
$ws.Range("M132").Value = -13334.8892
$ws.Range("N132").Value = -18179

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2719613.8
$ws.Range("I105").Value = 2843141.8
$ws.Range("K105").Value = 2843141.8
$ws.Range("M105").Value = -2841394.8

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 461.36365
$ws.Range("I6").Value = 492.5
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 492.5
$ws.Range("L6").Value = 150
$ws.Range("M6").Value = -379.5
$ws.Range("N6").Value = -376

# Row 33
$ws.Range("H33").Value = 2615.0908
$ws.Range("I33").Value = 1005.5
$ws.Range("J33").Value = 4546.6
$ws.Range("K33").Value = 1005.5
$ws.Range("L33").Value = 4546.6
$ws.Range("M33").Value = -626.5
$ws.Range("N33").Value = -5304.6

# Row 58
$ws.Range("H58").Value = 1829.4
$ws.Range("I58").Value = 1479.5333
$ws.Range("J58").Value = 2879
$ws.Range("K58").Value = 1479.5333
$ws.Range("L58").Value = 2879
$ws.Range("M58").Value = -1276.5333
$ws.Range("N58").Value = -3285

# Row 132
$ws.Range("H132").Value = 54939.473
$ws.Range("I132").Value = 57908.445
$ws.Range("J132").Value = 1498
$ws.Range("K132").Value = 173725.335
$ws.Range("L132").Value = 4494
$ws.Range("M132").Value = -171195.335
$ws.Range("N132").Value = -9554

# Row 136
$ws.Range("H136").Value = 1829.4
$ws.Range("I136").Value = 1479.5333
$ws.Range("J136").Value = 2879
$ws.Range("K136").Value = 4438.5999
$ws.Range("L136").Value = 8637
$ws.Range("M136").Value = -1888.5999
$ws.Range("N136").Value = -13737

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 441.5
$ws.Range("I34").Value = 183
$ws.Range("K34").Value = 549
$ws.Range("M34").Value = -465

# Row 46
$ws.Range("H46").Value = 168214.3
$ws.Range("I46").Value = 667398.2
$ws.Range("J46").Value = 1819.6666
$ws.Range("K46").Value = 2002194.6
$ws.Range("L46").Value = 5458.9998
$ws.Range("M46").Value = -2002103.6
$ws.Range("N46").Value = -5640.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Range("H44").Value = 26007
$ws.Range("I44").Value = 17014
$ws.Range("J44").Value = 35000
$ws.Range("K44").Value = 17014
$ws.Range("L44").Value = 35000
$ws.Range("M44").Value = -16418
$ws.Range("N44").Value = -36192

# Row 52
$ws.Range("H52").Value = 29673.2
$ws.Range("J52").Value = 29673.2
$ws.Range("L52").Value = 29673.2
$ws.Range("N52").Value = -30191.2

# Row 107
$ws.Range("H107").Value = 1266.7
$ws.Range("I107").Value = 1266.7
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1266.7
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 653.3
$ws.Range("N107").ClearContents()

# Row 108
$ws.Range("H108").Value = 100000
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2544.111
$ws.Range("I22").Value = 2400
$ws.Range("J22").Value = 2585.2856
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 2585.2856
$ws.Range("M22").Value = -2105
$ws.Range("N22").Value = -3175.2856

# Row 27
$ws.Range("H27").Value = 2544.111
$ws.Range("I27").Value = 2400
$ws.Range("J27").Value = 2585.2856
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 2585.2856
$ws.Range("M27").Value = -2293
$ws.Range("N27").Value = -2799.2856

# Row 46
$ws.Range("H46").Value = 1678717.1
$ws.Range("I46").Value = 2903802.8
$ws.Range("J46").Value = 8145.636
$ws.Range("K46").Value = 2903802.8
$ws.Range("L46").Value = 8145.636
$ws.Range("M46").Value = -2903614.8
$ws.Range("N46").Value = -8521.636

# Row 48
$ws.Range("H48").Value = 34583
$ws.Range("I48").Value = 34166
$ws.Range("K48").Value = 34166
$ws.Range("M48").Value = -33505

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 76
$ws.Range("H76").Value = 6765
$ws.Range("I76").Value = 261
$ws.Range("J76").Value = 8933
$ws.Range("K76").Value = 261
$ws.Range("L76").Value = 8933
$ws.Range("N76").Value = -9609
$ws.Range("M76").Value = 77

# Row 79
$ws.Range("H79").Value = 6765
$ws.Range("I79").Value = 261
$ws.Range("J79").Value = 8933
$ws.Range("K79").Value = 261
$ws.Range("L79").Value = 8933
$ws.Range("N79").Value = -11273
$ws.Range("M79").Value = 909

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 772.5217
$ws.Range("I113").Value = 350.46667
$ws.Range("K113").Value = 1051.40001
$ws.Range("M113").Value = 1118.59999

# Row 136
$ws.Range("H136").Value = 5652.654
$ws.Range("I136").Value = 7355.5293
$ws.Range("J136").Value = 2436.111
$ws.Range("K136").Value = 22066.5879
$ws.Range("L136").Value = 7308.333
$ws.Range("M136").Value = -19516.5879
